$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 19.9522100650596
$ws.Range("C2").Value = 9.523817997876256
$ws.Range("D2").Value = 7.263747435840701
$ws.Range("E2").Value = 9.389834134098166
$ws.Range("F2").Value = 39.11008160143846
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 30.67827950301364
$ws.Range("L2").Value = 10.49374847646194
# Row 3
$ws.Range("B3").Value = 19.51527105547705
$ws.Range("C3").Value = 8.959559772404385
$ws.Range("D3").Value = 7.282965234281217
$ws.Range("E3").Value = 9.405329733874636
$ws.Range("F3").Value = 38.73492891471717
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 30.58314328341842
$ws.Range("L3").Value = 10.47525806446455
# Row 4
$ws.Range("B4").Value = 19.24852251537796
$ws.Range("C4").Value = 8.596659747759007
$ws.Range("D4").Value = 7.295720102052073
$ws.Range("E4").Value = 9.415484258095537
$ws.Range("F4").Value = 38.5145948598505
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 30.53209033613327
$ws.Range("L4").Value = 10.46620087574232
# Row 5
$ws.Range("B5").Value = 19.14038131367918
$ws.Range("C5").Value = 8.444775081115568
$ws.Range("D5").Value = 7.301156872693566
$ws.Range("E5").Value = 9.419783543343964
$ws.Range("F5").Value = 38.42740137791341
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 30.51313955152016
$ws.Range("L5").Value = 10.46308911732519
# Row 6
$ws.Range("B6").Value = 19.12246361237931
$ws.Range("C6").Value = 8.419317679726648
$ws.Range("D6").Value = 7.302074050889154
$ws.Range("E6").Value = 9.420507182230978
$ws.Range("F6").Value = 38.41308179345118
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 30.51010478251477
$ws.Range("L6").Value = 10.46260743214762
# Row 7
$ws.Range("B7").Value = 19.24706158014746
$ws.Range("C7").Value = 8.594627374023233
$ws.Range("D7").Value = 7.295792457780266
$ws.Range("E7").Value = 9.415541586568592
$ws.Range("F7").Value = 38.51340833639572
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 30.53182725196273
$ws.Range("L7").Value = 10.46615656266549
# Row 8
$ws.Range("B8").Value = 19.80133929102343
$ws.Range("C8").Value = 9.332746818840914
$ws.Range("D8").Value = 7.270174846211089
$ws.Range("E8").Value = 9.395044307103618
$ws.Range("F8").Value = 38.97869918495624
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 30.64395101892312
$ws.Range("L8").Value = 10.48689764756647
# Row 9
$ws.Range("B9").Value = 20.89294613006417
$ws.Range("C9").Value = 10.70302579325677
$ws.Range("D9").Value = 7.227566252419511
$ws.Range("E9").Value = 9.359916952761882
$ws.Range("F9").Value = 39.9669030271473
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 30.9220183063048
$ws.Range("L9").Value = 10.54569146654217
# Row 10
$ws.Range("B10").Value = 21.68815800892231
$ws.Range("C10").Value = 11.65592536534008
$ws.Range("D10").Value = 7.200979944226457
$ws.Range("E10").Value = 9.337182297514707
$ws.Range("F10").Value = 40.73387671377589
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 31.16136813160077
$ws.Range("L10").Value = 10.59977624191547
# Row 11
$ws.Range("B11").Value = 22.04655271727565
$ws.Range("C11").Value = 12.06469726788564
$ws.Range("D11").Value = 7.18992462926238
$ws.Range("E11").Value = 9.327503741308254
$ws.Range("F11").Value = 41.09048443099082
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 31.27773105724922
$ws.Range("L11").Value = 10.62670200037033
# Row 12
$ws.Range("B12").Value = 22.18163975910087
$ws.Range("C12").Value = 12.21595684453953
$ws.Range("D12").Value = 7.185888941332425
$ws.Range("E12").Value = 9.32393389675793
$ws.Range("F12").Value = 41.22652624670572
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 31.32285553176661
$ws.Range("L12").Value = 10.63722759764472
# Row 13
$ws.Range("B13").Value = 22.15257657828222
$ws.Range("C13").Value = 12.1835372329176
$ws.Range("D13").Value = 7.186751374371807
$ws.Range("E13").Value = 9.32469849515142
$ws.Range("F13").Value = 41.1971843396047
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 31.31309027665774
$ws.Range("L13").Value = 10.63494614906975
# Row 14
$ws.Range("B14").Value = 22.0576798307479
$ws.Range("C14").Value = 12.07721218969023
$ws.Range("D14").Value = 7.189589582499603
$ws.Range("E14").Value = 9.327208140953198
$ws.Range("F14").Value = 41.10165720053406
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 31.28142231951292
$ws.Range("L14").Value = 10.62756137361978
# Row 15
$ws.Range("B15").Value = 21.99946667488261
$ws.Range("C15").Value = 12.01162544166568
$ws.Range("D15").Value = 7.191347734661709
$ws.Range("E15").Value = 9.328757765271874
$ws.Range("F15").Value = 41.04327141083486
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 31.26216237793171
$ws.Range("L15").Value = 10.6230807400358
# Row 16
$ws.Range("B16").Value = 21.66465582365625
$ws.Range("C16").Value = 11.62871564968221
$ws.Range("D16").Value = 7.20172345369761
$ws.Range("E16").Value = 9.337828148803782
$ws.Range("F16").Value = 40.71071762306261
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 31.15391281154286
$ws.Range("L16").Value = 10.59806293253547
# Row 17
$ws.Range("B17").Value = 21.4582940436476
$ws.Range("C17").Value = 11.38749865246183
$ws.Range("D17").Value = 7.208355652137791
$ws.Range("E17").Value = 9.343562340819176
$ws.Range("F17").Value = 40.50860579131532
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 31.08941104073567
$ws.Range("L17").Value = 10.58330702486891
# Row 18
$ws.Range("B18").Value = 21.33929300318175
$ws.Range("C18").Value = 11.24643033151882
$ws.Range("D18").Value = 7.212268026151149
$ws.Range("E18").Value = 9.346922964989735
$ws.Range("F18").Value = 40.39308836884135
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 31.05301658005004
$ws.Range("L18").Value = 10.57503870625486
# Row 19
$ws.Range("B19").Value = 21.29895317885883
$ws.Range("C19").Value = 11.19826682144454
$ws.Range("D19").Value = 7.21360943030642
$ws.Range("E19").Value = 9.348071549903054
$ws.Range("F19").Value = 40.35410507375902
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 31.04081560647282
$ws.Range("L19").Value = 10.57227692139437
# Row 20
$ws.Range("B20").Value = 21.48029449949741
$ws.Range("C20").Value = 11.41341727615783
$ws.Range("D20").Value = 7.207639521267476
$ws.Range("E20").Value = 9.342945462518353
$ws.Range("F20").Value = 40.53004594952854
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 31.0962044790875
$ws.Range("L20").Value = 10.5848551937851
# Row 21
$ws.Range("B21").Value = 22.0855715105438
$ws.Range("C21").Value = 12.10853817571899
$ws.Range("D21").Value = 7.18875183084786
$ws.Range("E21").Value = 9.326468414647168
$ws.Range("F21").Value = 41.12968949046304
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 31.29069532156301
$ws.Range("L21").Value = 10.62972155989268
# Row 22
$ws.Range("B22").Value = 22.47742571362571
$ws.Range("C22").Value = 12.54225302758589
$ws.Range("D22").Value = 7.177286770266905
$ws.Range("E22").Value = 9.316254604060225
$ws.Range("F22").Value = 41.52738273510347
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 31.42397834378222
$ws.Range("L22").Value = 10.66096188282135
# Row 23
$ws.Range("B23").Value = 22.26867317683506
$ws.Range("C23").Value = 12.31264856866924
$ws.Range("D23").Value = 7.183325006219819
$ws.Range("E23").Value = 9.321655200258956
$ws.Range("F23").Value = 41.31463185786355
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 31.35228357278137
$ws.Range("L23").Value = 10.64411449410309
# Row 24
$ws.Range("B24").Value = 21.47034920625877
$ws.Range("C24").Value = 11.4017069104727
$ws.Range("D24").Value = 7.207962974417146
$ws.Range("E24").Value = 9.343224153959056
$ws.Range("F24").Value = 40.52035072991106
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 31.09313102004068
$ws.Range("L24").Value = 10.58415459644806
# Row 25
$ws.Range("B25").Value = 20.59815977354553
$ws.Range("C25").Value = 10.33007424475663
$ws.Range("D25").Value = 7.238269434440181
$ws.Range("E25").Value = 9.368878968164305
$ws.Range("F25").Value = 39.6919812153371
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 30.84059858920514
$ws.Range("L25").Value = 10.52786066065583
